$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "AutoOppabc"
$ws.Range("A3").Value = "AutoOppxyz"
$ws.Range("A3").Select()
